$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update Coin name / Link (B, C) for rows whose ranking swapped ---
$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("B48").Value = 'Celestia'
$ws.Range("C48").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("B49").Value = 'SynthetixNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'

# --- Force Text number format on Price (D) cells whose new value looks numeric, ---
# --- so Excel keeps them as text (matching the source data's text-typed prices) ---
$forceTextCells = @("D5","D6","D7","D9","D10","D11","D12","D14","D16","D17","D21","D22","D23","D25","D27","D29","D30","D31","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D47","D48","D49")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Write new Price (D) values ---
$ws.Range("D2").Value = '41.035.27'
$ws.Range("D3").Value = '2.203.95'
$ws.Range("D5").Value = '242.67'
$ws.Range("D6").Value = '0.624'
$ws.Range("D7").Value = '70.28'
$ws.Range("D9").Value = '0.546'
$ws.Range("D10").Value = '36.71'
$ws.Range("D11").Value = '0.0951'
$ws.Range("D12").Value = '57.61'
$ws.Range("D14").Value = '6.63'
$ws.Range("D15").Value = '2.530.42'
$ws.Range("D16").Value = '14.69'
$ws.Range("D17").Value = '0.835'
$ws.Range("D18").Value = '2.204.52'
$ws.Range("D19").Value = '40.942.30'
$ws.Range("D20").Value = '0.0₃0943'
$ws.Range("D21").Value = '72.51'
$ws.Range("D22").Value = '6.06'
$ws.Range("D23").Value = '230.33'
$ws.Range("D25").Value = '0.999'
$ws.Range("D27").Value = '2.41'
$ws.Range("D29").Value = '9.70'
$ws.Range("D30").Value = '170.07'
$ws.Range("D31").Value = '20.29'
$ws.Range("D33").Value = '0.124'
$ws.Range("D34").Value = '0.0702'
$ws.Range("D35").Value = '5.14'
$ws.Range("D36").Value = '4.59'
$ws.Range("D37").Value = '3.84'
$ws.Range("D38").Value = '23.84'
$ws.Range("D40").Value = '0.0272'
$ws.Range("D41").Value = '5.78'
$ws.Range("D42").Value = '62.58'
$ws.Range("D43").Value = '4.84'
$ws.Range("D44").Value = '0.194'
$ws.Range("D45").Value = '8.58'
$ws.Range("D47").Value = '0.0989'
$ws.Range("D48").Value = '10.46'
$ws.Range("D49").Value = '4.44'

# --- Restore default formatting on the forced-text cells (keeps values as text) ---
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).ClearFormats()
}

# --- Write new Volume(1h) (E) values ---
$ws.Range("E2").Value = '  -6.80%  '
$ws.Range("E3").Value = '  -7.14%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("E6").Value = '  -7.61%  '
$ws.Range("E7").Value = '  -5.46%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -11.28%  '
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("E11").Value = '  -7.56%  '
$ws.Range("E12").Value = '  -4.01%  '
$ws.Range("E13").Value = '  -4.39%  '
$ws.Range("E14").Value = '  -9.69%  '
$ws.Range("E15").Value = '  -7.26%  '
$ws.Range("E16").Value = '  -10.67%  '
$ws.Range("E17").Value = '  -10.23%  '
$ws.Range("E18").Value = '  -7.57%  '
$ws.Range("E19").Value = '  -6.84%  '
$ws.Range("E20").Value = '  -8.91%  '
$ws.Range("E21").Value = '  -7.15%  '
$ws.Range("E22").Value = '  -8.37%  '
$ws.Range("E23").Value = '  -9.69%  '
$ws.Range("E24").Value = '  +6.71%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -4.92%  '
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("E28").Value = '  -5.17%  '
$ws.Range("E29").Value = '  -8.51%  '
$ws.Range("E30").Value = '  -3.37%  '
$ws.Range("E31").Value = '  -9.62%  '
$ws.Range("E32").Value = '  -9.51%  '
$ws.Range("E33").Value = '  -8.15%  '
$ws.Range("E34").Value = '  -7.46%  '
$ws.Range("E35").Value = '  -5.05%  '
$ws.Range("E36").Value = '  -10.15%  '
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("E38").Value = '  +14.50%  '
$ws.Range("E39").Value = '  -6.52%  '
$ws.Range("E40").Value = '  -3.58%  '
$ws.Range("E41").Value = '  -13.19%  '
$ws.Range("E42").Value = '  -5.02%  '
$ws.Range("E43").Value = '  -12.60%  '
$ws.Range("E44").Value = '  -4.77%  '
$ws.Range("E45").Value = '  -5.85%  '
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("E47").Value = '  -8.36%  '
$ws.Range("E48").Value = '  +7.96%  '
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("E50").Value = '  -6.63%  '
$ws.Range("E51").Value = '  -6.53%  '
